$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.616.89'
$ws.Range("E2").Value = '  +1.13%  '

$ws.Range("D3").Value = '3.393.05'
$ws.Range("E3").Value = '  -0.02%  '

$ws.Range("E4").Value = '  -0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '575.75'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.61%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '141.24'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.76%  '

$ws.Range("E8").Value = '  -0.56%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.70'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.77%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.123'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.02%  '

$ws.Range("E11").Value = '  -2.48%  '

$ws.Range("D12").Value = '3.967.93'
$ws.Range("E12").Value = '  -0.07%  '

$ws.Range("E13").Value = '  +0.16%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.40'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.27%  '

$ws.Range("D15").Value = '3.391.54'
$ws.Range("E15").Value = '  +0.11%  '

$ws.Range("E16").Value = '  -0.22%  '

$ws.Range("D17").Value = '61.643.52'
$ws.Range("E17").Value = '  +1.06%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.13'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.18%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.64'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.41%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '9.03'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.76%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '391.48'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.18%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '74.98'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.93%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.551'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.32%  '

$ws.Range("E24").Value = '  +0.26%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000113'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.73%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.194'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +7.89%  '

$ws.Range("E27").Value = '  -0.06%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.28'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.27%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.05'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.31%  '

$ws.Range("E30").Value = '  -0.81%  '

$ws.Range("E31").Value = '  +0.00%  '

$ws.Range("E32").Value = '  +0.05%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '23.30'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.76%  '

$ws.Range("B34").Value = 'Aptos'
$ws.Range("C34").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.92'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.17%  '

$ws.Range("B35").Value = 'Monero'
$ws.Range("C35").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '169.76'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.19%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.05'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.75%  '

$ws.Range("D37").Value = '3.424.93'
$ws.Range("E37").Value = '  +0.03%  '

$ws.Range("E38").Value = '  -0.83%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0768'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.48%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '26.14'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -4.77%  '

$ws.Range("E41").Value = '  -0.14%  '

$ws.Range("E42").Value = '  -0.16%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.66'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.78%  '

$ws.Range("E44").Value = '  +2.31%  '

$ws.Range("D45").Value = '2.471.80'
$ws.Range("E45").Value = '  -0.54%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '22.87'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.64%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '6.67'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.96%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.00'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.03%  '

$ws.Range("E49").Value = '  -0.89%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.03'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.59%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.207'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.40%  '
